$d = $word.ActiveDocument

# Locate the literal run of text "{COMPANY_NAME}" (the brace-delimited
# placeholder immediately after "offered by ") and collapse the found
# range to its start so we can insert a new "$" character immediately
# before the opening brace, matching the author's edit which prefixed
# the placeholder with a dollar sign (as already done for {JOB_TITAL}).
$rng = $d.Content
$found = $rng.Find.Execute("{COMPANY_NAME}", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

$rng.Collapse(1)
$rng.InsertBefore("$")
